$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H13").Value = 34000
$ws.Range("I13").Value = 34000
$ws.Range("J13").Value = 0
$ws.Range("K13").Value = 34000
$ws.Range("L13").Value = 0
$ws.Range("M13").ClearContents()
$ws.Range("N13").ClearContents()
$ws.Range("H28").Value = 843.6
$ws.Range("I28").Value = 1267.3077
$ws.Range("J28").Value = 56.714287
$ws.Range("K28").Value = 1267.3077
$ws.Range("L28").Value = 56.714287
$ws.Range("M28").Value = -782.3077000000001
$ws.Range("N28").Value = -1026.714287
$ws.Range("H113").Value = 3201.037
$ws.Range("J113").Value = 2867.6924
$ws.Range("L113").Value = 2867.6924
$ws.Range("N113").Value = -9375.6924

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 976.48
$ws.Range("I32").Value = 976.48
$ws.Range("J32").Value = 0
$ws.Range("K32").Value = 976.48
$ws.Range("L32").Value = 0
$ws.Range("M32").ClearContents()
$ws.Range("N32").ClearContents()
$ws.Range("H45").Value = 1411.68
$ws.Range("I45").Value = 978
$ws.Range("J45").Value = 2785
$ws.Range("K45").Value = 978
$ws.Range("L45").Value = 2785
$ws.Range("M45").Value = -601
$ws.Range("N45").Value = -3539
$ws.Range("H61").Value = 1181.4067
$ws.Range("I61").Value = 1211.0204
$ws.Range("K61").Value = 1211.0204
$ws.Range("M61").Value = -999.0204000000001
$ws.Range("H132").Value = 1121.6719
$ws.Range("I132").Value = 928.56366
$ws.Range("J132").Value = 2301.7778
$ws.Range("K132").Value = 2785.69098
$ws.Range("L132").Value = 6905.3334
$ws.Range("M132").Value = -255.6909800000003
$ws.Range("N132").Value = -11965.3334
$ws.Range("H136").Value = 1181.4067
$ws.Range("I136").Value = 1211.0204
$ws.Range("K136").Value = 3633.0612
$ws.Range("M136").Value = -1083.0612

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H16").Value = 35000
$ws.Range("J16").Value = 35000
$ws.Range("L16").Value = 35000
$ws.Range("N16").Value = -35340
$ws.Range("H134").Value = 24764.75
$ws.Range("I134").Value = 1581.0526
$ws.Range("K134").Value = 4743.1578
$ws.Range("M134").Value = -2208.1578

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1114.2106
$ws.Range("I16").Value = 996
$ws.Range("J16").Value = 1276.75
$ws.Range("K16").Value = 996
$ws.Range("L16").Value = 1276.75
$ws.Range("M16").Value = -709
$ws.Range("N16").Value = -1850.75
$ws.Range("H58").Value = 3260.4695
$ws.Range("I58").Value = 1119.425
$ws.Range("J58").Value = 12776.223
$ws.Range("K58").Value = 1119.425
$ws.Range("L58").Value = 12776.223
$ws.Range("M58").Value = -916.425
$ws.Range("N58").Value = -13182.223
$ws.Range("H113").Value = 1114.2106
$ws.Range("I113").Value = 996
$ws.Range("J113").Value = 1276.75
$ws.Range("K113").Value = 996
$ws.Range("L113").Value = 1276.75
$ws.Range("M113").Value = 1174
$ws.Range("N113").Value = -5616.75
$ws.Range("H134").Value = 14286812
$ws.Range("I134").Value = 1044.5667
$ws.Range("J134").Value = 100001416
$ws.Range("K134").Value = 3133.7001
$ws.Range("L134").Value = 300004248
$ws.Range("M134").Value = -598.7001
$ws.Range("N134").Value = -300009318
$ws.Range("H136").Value = 3260.4695
$ws.Range("I136").Value = 1119.425
$ws.Range("J136").Value = 12776.223
$ws.Range("K136").Value = 3358.275
$ws.Range("L136").Value = 38328.669
$ws.Range("M136").Value = -808.2749999999996
$ws.Range("N136").Value = -43428.669

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H26").Value = 265
$ws.Range("J26").Value = 298
$ws.Range("L26").Value = 894
$ws.Range("N26").Value = -1470
$ws.Range("H108").Value = 333
$ws.Range("I108").Value = 333
$ws.Range("J108").Value = 0
$ws.Range("K108").Value = 999
$ws.Range("L108").Value = 0
$ws.Range("M108").ClearContents()
$ws.Range("N108").ClearContents()
$ws.Range("H110").Value = 11276.637
$ws.Range("I110").Value = 4510.75
$ws.Range("J110").Value = 15142.857
$ws.Range("K110").Value = 13532.25
$ws.Range("L110").Value = 45428.571
$ws.Range("M110").Value = -9442.25
$ws.Range("N110").Value = -53608.571
$ws.Range("H111").Value = 0
$ws.Range("I111").Value = 0
$ws.Range("J111").Value = 0
$ws.Range("K111").Value = 0
$ws.Range("L111").ClearContents()
$ws.Range("M111").ClearContents()
$ws.Range("N111").ClearContents()
$ws.Range("H112").Value = 4003.375
$ws.Range("I112").Value = 2513.5
$ws.Range("J112").Value = 4500
$ws.Range("K112").Value = 7540.5
$ws.Range("L112").Value = 13500
$ws.Range("M112").Value = -6432.5
$ws.Range("N112").Value = -15716
$ws.Range("H113").Value = 874.61536
$ws.Range("I113").Value = 1532.6666
$ws.Range("J113").Value = 677.2
$ws.Range("K113").Value = 4597.9998
$ws.Range("L113").Value = 2031.6
$ws.Range("M113").Value = -2427.9998
$ws.Range("N113").Value = -6371.6
$ws.Range("H114").Value = 1442
$ws.Range("I114").Value = 664.6667
$ws.Range("J114").Value = 2996.6667
$ws.Range("K114").Value = 1994.0001
$ws.Range("L114").Value = 8990.000100000001
$ws.Range("M114").Value = 1259.9999
$ws.Range("N114").Value = -15498.0001
$ws.Range("H115").Value = 1391
$ws.Range("I115").Value = 982.3333
$ws.Range("K115").Value = 2946.9999
$ws.Range("M115").Value = -1771.9999
$ws.Range("H116").Value = 168800
$ws.Range("I116").Value = 2560
$ws.Range("J116").Value = 1000000
$ws.Range("K116").Value = 7680
$ws.Range("L116").Value = 3000000
$ws.Range("M116").Value = -4238
$ws.Range("N116").Value = -3006884
$ws.Range("H117").Value = 2584.4666
$ws.Range("I117").Value = 476.33334
$ws.Range("J117").Value = 3111.5
$ws.Range("K117").Value = 1429.00002
$ws.Range("L117").Value = 9334.5
$ws.Range("M117").Value = 2012.99998
$ws.Range("N117").Value = -16218.5
$ws.Range("H118").Value = 2728
$ws.Range("I118").Value = 400
$ws.Range("J118").Value = 3504
$ws.Range("K118").Value = 1200
$ws.Range("L118").Value = 10512
$ws.Range("M118").Value = 43
$ws.Range("N118").Value = -12998
$ws.Range("H119").Value = 13428.429
$ws.Range("I119").Value = 8166.3335
$ws.Range("J119").Value = 17375
$ws.Range("K119").Value = 24499.0005
$ws.Range("L119").Value = 52125
$ws.Range("M119").Value = -19661.0005
$ws.Range("N119").Value = -61801
$ws.Range("H120").Value = 7156900
$ws.Range("I120").Value = 50007500
$ws.Range("J120").Value = 15133.333
$ws.Range("K120").Value = 150022500
$ws.Range("L120").Value = 45399.999
$ws.Range("M120").Value = -150017662
$ws.Range("N120").Value = -55075.999
$ws.Range("H121").Value = 956.6923
$ws.Range("J121").Value = 956.6923
$ws.Range("L121").Value = 2870.0769
$ws.Range("N121").Value = -5490.0769
$ws.Range("H131").Value = 17277518
$ws.Range("J131").Value = 1551.2084
$ws.Range("L131").Value = 4653.6252
$ws.Range("N131").Value = -14733.6252

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 4738.231
$ws.Range("I70").Value = 4699.625
$ws.Range("J70").Value = 4800
$ws.Range("K70").Value = 4699.625
$ws.Range("L70").Value = 4800
$ws.Range("M70").Value = -4429.625
$ws.Range("N70").Value = -5340
$ws.Range("H73").Value = 4738.231
$ws.Range("I73").Value = 4699.625
$ws.Range("J73").Value = 4800
$ws.Range("K73").Value = 4699.625
$ws.Range("L73").Value = 4800
$ws.Range("M73").Value = -3763.625
$ws.Range("N73").Value = -6672
$ws.Range("H132").Value = 3309.8462
$ws.Range("I132").Value = 3018.8948
$ws.Range("J132").Value = 4099.5713
$ws.Range("K132").Value = 9056.6844
$ws.Range("L132").Value = 12298.7139
$ws.Range("M132").Value = -6526.6844
$ws.Range("N132").Value = -17358.7139

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H110").Value = 36483.332
$ws.Range("J110").Value = 36483.332
$ws.Range("L110").Value = 36483.332
$ws.Range("N110").Value = -44663.332
